# Online Enrollment Check List - add two new TABLES notes under the
# "Insert Membership information into Production tables" block (D24/D25),
# pushing the existing row 25 and everything below it down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 25 (old row 25 and everything below shifts down by one).
$ws.Rows(25).Insert()

# Fill the freshly exposed D24/D25 cells with the new check-list notes.
$ws.Range("D24").Value = "Taxes in asprecdoc? Seems mixed on existing - ask John/Mel"
$ws.Range("D25").Value = "offer PIN/Password to customer"

# Match the author's final selection/active cell.
$ws.Range("D25").Select() | Out-Null
